$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date value for rows 2-15 from 45174 to 45175
$oldDate = [DateTime]::FromOADate(45174)
for ($r = 2; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value() -eq $oldDate) {
        $cell.Value = 45175
    }
}
